$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)
$sh = $s.Shapes.Item(1)
# Target position is -2909, 137922 EMU (1 pt = 12700 EMU).
# PowerPoint's Shape.Top is a single-precision float, and the EMU value it
# round-trips to is truncated rather than rounded; nudge slightly above the
# exact quotient so the stored single still truncates to the exact EMU target.
$sh.Top = 10.8600005
